$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Device")

# Rename the "Device" sheet to "Apparatus"
$ws.Name = "Apparatus"

# Update the sheet's descriptive note to refer to "apparatuses" instead of "devices"
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."

# Make the header row (bus number / type number / parameters) bold, matching
# the formatting already used for the title/"Data:" rows above it
$ws.Range("A1").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)

# Reflect the newly-selected header range
$ws.Range("A3:C3").Select() | Out-Null
